$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'53.854.17"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -11.14%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.310.35"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -20.43%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  -0.14%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'443.80"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -15.90%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'126.09"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -12.37%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.996"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -0.38%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.474"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -14.93%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'2.304.19"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -20.88%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'5.36"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -11.12%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  -15.40%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.309"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -14.64%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "'  -3.57%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'2.709.33"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -20.57%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'53.816.87"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -11.22%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'18.76"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -17.53%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("E17").Value = "'  -14.68%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'2.320.59"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -20.30%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'4.01"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -20.78%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'297.87"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -17.64%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'9.32"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -20.37%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'  -0.19%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'5.58"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -1.98%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'5.36"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -19.03%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'55.36"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -14.62%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +0.40%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'0.154"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -15.06%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'  -19.11%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'6.89"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -12.48%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'0.998"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -0.14%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'0.0₃0704"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -18.13%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'146.36"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -3.76%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'16.88"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -14.59%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  -19.86%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'4.67"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -16.24%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'3.57"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -18.73%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'0.835"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -17.49%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  -16.50%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.993"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -0.40%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'33.17"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -12.13%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  -0.28%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'  -15.82%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'  -18.05%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'1.922.56"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -16.27%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = "'  -15.28%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.520"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -19.68%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.0209"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -12.10%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.0825"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -10.78%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'15.98"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -22.41%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'3.99"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -20.22%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'4.65"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -3.63%  "
$ws.Range("E51").Style = "Normal"
